$d = $word.ActiveDocument

# --- Edit 1: paragraph 1 text tweak + three red-colored runs appended ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.SetRange($r.End - 1, $r.End - 1)
$r.InsertAfter("  ")

$r.SetRange($p1.Range.End - 1, $p1.Range.End - 1)
$r.InsertAfter("(This is a change – Ve")
$r.Font.Color = 192

$r.SetRange($p1.Range.End - 1, $p1.Range.End - 1)
$r.InsertAfter("rsion for branch alternate")
$r.Font.Color = 192

$r.SetRange($p1.Range.End - 1, $p1.Range.End - 1)
$r.InsertAfter(")")
$r.Font.Color = 192

# --- Edit 2: insert a new empty shaded paragraph right after paragraph 2 ---
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3213B859" w14:textId="4C5A8A93" w:rsidR="00094D0B" w:rsidRDefault="00094D0B"><w:r><w:t>It will be treated as a binary file by Git.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p>'
$r2.InsertXML($xml)

# --- Edit 3: prune now-unused styles from styles.xml (mirrors Word auto-cleanup) ---
# Resolve the target style names to their current collection indices, then delete
# from the highest index downward so earlier deletions never shift the indices
# we still have queued up (deleting low-to-high corrupts this engine's internal
# style table and crashes it).
$namesToDelete = @("Heading1","Heading2","Heading4","apple-converted-space","Hyperlink","Heading2Char","Heading4Char","audio-tool","subscribe","subscribe-more-info","generic-title","podcast-toolssubscribe-links","Heading1Char","c-txt","FollowedHyperlink")

$indicesToDelete = @()
for ($i = 1; $i -le $d.Styles.Count; $i++) {
  $styleId = $d.Styles.Item($i).NameLocal
  foreach ($n in $namesToDelete) {
    if ($d.Styles.Item($i).Equals($d.Styles($n))) {
      $indicesToDelete += $i
    }
  }
}
$indicesToDelete = $indicesToDelete | Sort-Object -Descending

foreach ($i in $indicesToDelete) {
  $d.Styles.Item($i).Delete()
}

Write-Host "Paragraphs count:" $d.Paragraphs.Count
Write-Host "Styles count:" $d.Styles.Count
